$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 39 data: Problem 39 Java Backtracking 100% 71%
$ws.Range("A39").NumberFormat = "0.00%"
$ws.Range("A39").Value = 1

$ws.Range("B39").NumberFormat = "0%"
$ws.Range("B39").Value = 0.71

$ws.Range("C39").Value = "Backtracking"

# Update view to match: topLeftCell A31, selection C39
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C39").Select()
